# Change output caps for the 3.3V regulator section (row 6):
#   Old: A6 = "47uF ceramic"  (no part number, no qty)
#   New: A6 = "220uF electrolytic", B6 = "732-8911-1-ND", C6 = 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "732-8911-1-ND"
$ws.Range("A6").Value = "220uF electrolytic"
$ws.Range("C6").Value = 1

# Move the active selection to A6, matching the saved cursor position.
$ws.Range("A6").Select()
